$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0006408296065709695
$ws.Range("C2").Value = 0.04071648406533734
$ws.Range("D2").Value = 0.1494219747398047
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 0.6850158244724827

$ws.Range("B3").Value = 0.1190320826869504
$ws.Range("C3").Value = 0.04071648406533734
$ws.Range("D3").Value = 0.1494219747398047
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 0.8034070775528621
